$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the trailing "closing" row (26) so its leftover row-level
# formatting (thick bottom border / custom height) does not survive; we
# rebuild rows 26-40 below with the same plain style as rows 14-25. ---
$ws.Rows.Item(26).Delete()

# Copy the clean formatting of row 25 (fill/border only, no thick-bottom)
# down across the new rows 26-40.
$ws.Range("A25:E25").Copy()
$ws.Range("A26:E40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 20 never had an A-column cell before; give it the same look as its
# neighbours (A14:A19/A21:A25) by copying A19's format onto it.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the requirement text for column A, rows 13-40.
$arr = New-Object 'object[,]' 28,1
$arr[0,0] = "Terminal mobile de atendimento;"
$arr[1,0] = "Sensor de vagas livres/ocupadas"
$arr[2,0] = "Controle de filas de entrada e saída"
$arr[3,0] = "Controle de tickets"
$arr[4,0] = "Tratamento de meios de pagamento"
$arr[5,0] = "Tipificação de uso (mensalista/avulso/conveniado)"
$arr[6,0] = "Administração de caixa"
$arr[7,0] = "Configuração de tabelas preços (dias/horários/tipo de uso)"
$arr[8,0] = "Emissão de nota fiscal"
$arr[9,0] = "Gestão de cobranças de mensalistas e conveniados"
$arr[10,0] = "Reconhecimento de placa para recuperar dados cadastrais"
$arr[11,0] = "Cadastro de veículos e clientes"
$arr[12,0] = "Reservar vagas"
$arr[13,0] = "Solução mobile para o cliente fazer reservas e pagamentos"
$arr[14,0] = "Orientação por voz"
$arr[15,0] = "Integração com TAG"
$arr[16,0] = "Relatório de gerenciamento de média de ocupação de vagas por dia e horário"
$arr[17,0] = "Relatório de gerenciamento do valor recebido por meio de pagamento por mês"
$arr[18,0] = "Forum/Sistema de Denuncias"
$arr[19,0] = "Estacionamento Vertical com elevador"
$arr[20,0] = "Planta digital"
$arr[21,0] = "Sistema de fidelidade"
$arr[22,0] = "Controle do período estacionado"
$arr[23,0] = "Sistema de vigilância automatizado, com câmeras inteligentes. "
$arr[24,0] = "Ticket digital via e-mail ou sms"
$arr[25,0] = "Sistema de Segurança contra roubos "
$arr[26,0] = "Mostrar quantidade de vagas disponíveis em determinado local"
$arr[27,0] = "Sistema de redirecionamento para o condutor indicando vagas livres"
$ws.Range("A13:A40").Value = $arr

# Column widths tweak (A slightly wider, E slightly narrower).
$ws.Columns.Item(1).ColumnWidth = 74.7109375
$ws.Columns.Item(5).ColumnWidth = 86.85546875

# View state: scroll so row 28 is at the top, zoom to 85%, select E53.
$excel.ActiveWindow.Zoom = 85
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E53").Select()

Write-Host "Done"
